# Add a "Normal" (노멀) translation row to the translations sheet.
#
# In the source edit, a new row was inserted right above the existing
# row 22 ("Rare" grade translation), pushing rows 22:138 down to 23:139.
# The new row reuses the same cell formatting (style index 1 / text
# format) as its neighbouring rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 22; this shifts rows 22:138 down to
# 23:139 and carries the existing formatting onto the new row's cells,
# matching the surrounding rows.
$ws.Rows("22:22").Insert()

# Populate the newly inserted row with the Korean source text and its
# English translation.
$ws.Range("A22").Value = "노멀"
$ws.Range("B22").Value = "Normal"

# The workbook-level defined name "translations_1" describes the used
# range of the translation table; it must grow by one row to keep
# covering the whole table (was Book1!$A$1:$B$102, now $B$103).
$wb.Names("translations_1").RefersTo = "=Book1!`$A`$1:`$B`$103"

# Reflect the selection left behind by the edit (cell B22, the newly
# added translation) and scroll the view back near the top of the sheet.
$ws.Range("B22").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
